$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")

# Insert a new row at row 50 (shifts rows 50+ down by one), duplicating the
# "Hardware for ms41 case" spacer row that currently sits at row 49 -
# matching the author's "select row 49, copy, insert copied cells at row 50" edit.
$ws.Rows.Item(50).Insert(-4121)

# Copy formatting (styles/borders/number formats) from row 49 onto the newly
# inserted row 50 so the new row matches the template "spacer" row look.
$ws.Range("A49:S49").Copy()
$ws.Range("A50:S50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Bring across the label text that row 49 carries ("Hardware for ms41 case").
$ws.Range("B50").Value2 = $ws.Range("B49").Value2

# Re-create the two helper formulas (P and R columns) that every spacer /
# data row in this table carries, anchored to row 50.
$ws.Range("P50").Formula = '=IF(NOT(I50=""),A50&","&I50,"")'
$ws.Range("R50").Formula = '=IF(NOT(J50=""),J50&"|"&A50,"")'

# The row-wide format paste also stamped Q50/S50 (row 49 has no content
# there either) - make sure they stay empty like the rest of the spacer rows.
$ws.Range("Q50").Clear()
$ws.Range("S50").Clear()

# Match the author's final on-screen selection: the whole new row 50.
$ws.Rows.Item(50).Select()
